# Update "想去人数" (number of people interested) counts that changed
# between the two data pulls, on both the "展览" sheet and the
# "全部类型" sheet (which mirrors the same events).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 282
    $ws.Range("F3").Value = 16
    $ws.Range("F4").Value = 7768
    $ws.Range("F5").Value = 5659
    $ws.Range("F6").Value = 471
    $ws.Range("F7").Value = 76
}

# "展览" sheet has this event on row 11, while "全部类型" has it on row 13
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F11").Value = 246

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F13").Value = 246
